# Oct 18th notes update — append the Maven "pom.xml" / "Parent-Child" / "Release life cycle"
# sections to the bottom of the "Build Tools" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Build Tools")

function Set-BlankCell($addr) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Helvetica"
    $c.Font.Size = 12
}

function Set-TextCell($addr, $text) {
    $c = $ws.Range($addr)
    $c.Value = $text
    $c.Font.Name = "Helvetica Neue"
    $c.Font.Size = 16
    $c.Font.Color = 0
}

# ---- Row 79: "pom.xml" header -----------------------------------------
Set-TextCell "C79" "pom.xml"
Set-BlankCell "D79"
$ws.Rows.Item(79).RowHeight = 20

# ---- Rows 80-90: GAV / pom.xml element notes --------------------------
Set-BlankCell "C80"
Set-TextCell "D80" "<xml> version and opening of project"
$ws.Rows.Item(80).RowHeight = 20

Set-BlankCell "C81"
Set-TextCell "D81" "GAV"
$ws.Rows.Item(81).RowHeight = 20

Set-BlankCell "C82"
Set-TextCell "D82" "packaging —> jar/war/tar"
$ws.Rows.Item(82).RowHeight = 20

Set-BlankCell "C83"
Set-TextCell "D83" "dependencies —> we define all our project dependencies"
$ws.Rows.Item(83).RowHeight = 20

Set-BlankCell "C84"
Set-TextCell "D84" "modules —> To define all child components."
$ws.Rows.Item(84).RowHeight = 20

Set-BlankCell "C85"
Set-TextCell "D85" "dependencyManagement"
$ws.Rows.Item(85).RowHeight = 20

Set-BlankCell "C86"
Set-TextCell "D86" "scm —> we define our SCM URL"
$ws.Rows.Item(86).RowHeight = 20

Set-BlankCell "C87"
Set-TextCell "D87" "repositories —> artifact Repository URL"
$ws.Rows.Item(87).RowHeight = 20

Set-BlankCell "C88"
Set-TextCell "D88" "pluginRepositories —> Apache Maven’s Repo URL"
$ws.Rows.Item(88).RowHeight = 20

Set-BlankCell "C89"
Set-TextCell "D89" "distributionManagement —> Artifact repository URL and will be used in “mvn deploy” phase"
$ws.Rows.Item(89).RowHeight = 20

Set-BlankCell "C90"
Set-TextCell "D90" "profiles —> based on project/user/…"
$ws.Rows.Item(90).RowHeight = 20

# ---- Rows 91-92: blank spacer rows ------------------------------------
Set-BlankCell "C91"
Set-BlankCell "D91"

Set-BlankCell "C92"
Set-BlankCell "D92"

# ---- Row 93: "Parent/Child Module" header -----------------------------
Set-TextCell "C93" "Parent/Child Module"
Set-BlankCell "D93"
$ws.Rows.Item(93).RowHeight = 20

# ---- Row 94: "October" -> "pom.xml -> parent POM" ---------------------
Set-TextCell "C94" "October"
Set-TextCell "D94" "pom.xml -> parent POM"
$ws.Rows.Item(94).RowHeight = 20

# ---- Rows 95-98: Wk1..Wk4 pom.xml --------------------------------------
Set-BlankCell "C95"
Set-TextCell "D95" "Wk1/pom.xml"
$ws.Rows.Item(95).RowHeight = 20

Set-BlankCell "C96"
Set-TextCell "D96" "Wk2/pom.xml"
$ws.Rows.Item(96).RowHeight = 20

Set-BlankCell "C97"
Set-TextCell "D97" "Wk3/pom.xml"
$ws.Rows.Item(97).RowHeight = 20

Set-BlankCell "C98"
Set-TextCell "D98" "Wk4/pom.xml"
$ws.Rows.Item(98).RowHeight = 20

# ---- Rows 99-101: blank spacer rows ------------------------------------
Set-BlankCell "C99"
Set-BlankCell "D99"

Set-BlankCell "C100"
Set-BlankCell "D100"

Set-BlankCell "C101"
Set-BlankCell "D101"

# ---- Row 102: "Release life cycle" header ------------------------------
Set-TextCell "C102" "Release life cycle"
Set-BlankCell "D102"
$ws.Rows.Item(102).RowHeight = 20

# ---- Rows 103-105: Prepare / Perform / Rollback ------------------------
Set-BlankCell "C103"
Set-TextCell "D103" "Prepare —> mvn release:prepare"
$ws.Rows.Item(103).RowHeight = 20

Set-BlankCell "C104"
Set-TextCell "D104" "Perform —> mvn release:perform ( 1.0-SNPASHOT -> 1.0, 1.1/2.0-SNAPSHOT)"
$ws.Rows.Item(104).RowHeight = 20

Set-BlankCell "C105"
Set-TextCell "D105" "Rollback —> mvn release:rollback"
$ws.Rows.Item(105).RowHeight = 20

# Merge the "Release life cycle" label's blank C column alongside its 3 detail rows
$ws.Range("C103:C105").Merge()

# Scroll / select near the newly-added content, like the author did while typing it
$ws.Application.ActiveWindow.ScrollRow = 81
$ws.Range("D99").Select()
